# sd/qa/unit/data/pptx/smartart-linear-rule.pptx
#
# The underlying change moves the SmartArt diagram's graphic frame
# further down the slide (its vertical offset changes from
# 1407600 EMU to 2847600 EMU, i.e. +1440000 EMU == +1.5cm / +113.39pt),
# while everything else about the frame (width/height) stays the same.
#
# PowerPoint expresses shape position in points (1 pt == 12700 EMU), so
# convert the EMU offsets accordingly.

$EMU_PER_POINT = 12700
$newTopEmu = 2847600

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$diagram = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasSmartArt -or $shp.Name -eq "Diagram1") {
        $diagram = $shp
        break
    }
}
if ($null -eq $diagram) {
    $diagram = $s.Shapes.Item(1)
}

$diagram.Top = $newTopEmu / $EMU_PER_POINT
